$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 6667
$ws.Range("I12").Value = 7600
$ws.Range("K12").Value = 7600
$ws.Range("M12").Value = -7430
$ws.Range("H19").Value = 1024.75
$ws.Range("I19").Value = 1099.5
$ws.Range("K19").Value = 1099.5
$ws.Range("M19").Value = -924.5
$ws.Range("H55").Value = 246.25
$ws.Range("I55").Value = 296
$ws.Range("K55").Value = 296
$ws.Range("M55").Value = -82
$ws.Range("H62").Value = 4769.4443
$ws.Range("I62").Value = 4277.3076
$ws.Range("J62").Value = 6049
$ws.Range("K62").Value = 4277.3076
$ws.Range("L62").Value = 6049
$ws.Range("M62").Value = -3653.3076
$ws.Range("N62").Value = -7297
$ws.Range("H65").Value = 4769.4443
$ws.Range("I65").Value = 4277.3076
$ws.Range("J65").Value = 6049
$ws.Range("K65").Value = 21386.538
$ws.Range("L65").Value = 30245
$ws.Range("M65").Value = -18266.538
$ws.Range("N65").Value = -36485
$ws.Range("H86").Value = 1504.1
$ws.Range("I86").Value = 1469.75
$ws.Range("J86").Value = 1527
$ws.Range("K86").Value = 1469.75
$ws.Range("L86").Value = 1527
$ws.Range("M86").Value = -346.75
$ws.Range("N86").Value = -3773
$ws.Range("H89").Value = 1504.1
$ws.Range("I89").Value = 1469.75
$ws.Range("J89").Value = 1527
$ws.Range("K89").Value = 7348.75
$ws.Range("L89").Value = 7635
$ws.Range("M89").Value = -1732.75
$ws.Range("N89").Value = -18867
$ws.Range("H137").Value = 2286.524
$ws.Range("I137").Value = 2015.3572
$ws.Range("K137").Value = 6046.071599999999
$ws.Range("M137").Value = -3496.071599999999
$ws.Range("H138").Value = 2406.7273
$ws.Range("J138").Value = 2331.2273
$ws.Range("L138").Value = 6993.6819
$ws.Range("N138").Value = -17273.6819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3009.85
$ws.Range("I32").Value = 3037.5264
$ws.Range("K32").Value = 3037.5264
$ws.Range("M32").Value = -2750.5264
$ws.Range("H61").Value = 3429.5715
$ws.Range("I61").Value = 3167.8333
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 3167.8333
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2955.8333
$ws.Range("N61").Value = -5424
$ws.Range("H110").Value = 2501.25
$ws.Range("I110").Value = 2501.25
$ws.Range("K110").Value = 2501.25
$ws.Range("M110").Value = -456.25
$ws.Range("H136").Value = 3429.5715
$ws.Range("I136").Value = 3167.8333
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 9503.499899999999
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -6953.499899999999
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1376.2858
$ws.Range("I22").Value = 1094.8182
$ws.Range("K22").Value = 1094.8182
$ws.Range("M22").Value = -921.8181999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1839.3939
$ws.Range("I107").Value = 1507.5927
$ws.Range("K107").Value = 1507.5927
$ws.Range("M107").Value = 412.4073000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 20553666
$ws.Range("J4").Value = 81102290
$ws.Range("L4").Value = 243306870
$ws.Range("N4").Value = -243307094
$ws.Range("H54").Value = 10000
$ws.Range("J54").Value = 10000
$ws.Range("L54").Value = 30000
$ws.Range("N54").Value = -31118
$ws.Range("H58").Value = 7751.5
$ws.Range("I58").Value = 8000
$ws.Range("J58").Value = 7503
$ws.Range("K58").Value = 24000
$ws.Range("L58").Value = 22509
$ws.Range("M58").Value = -23872
$ws.Range("N58").Value = -22765
$ws.Range("H81").Value = 100013530
$ws.Range("I81").Value = 5875.3335
$ws.Range("J81").Value = 250025000
$ws.Range("K81").Value = 17626.0005
$ws.Range("L81").Value = 750075000
$ws.Range("M81").Value = -16503.0005
$ws.Range("N81").Value = -750077246
$ws.Range("H84").Value = 100013530
$ws.Range("I84").Value = 5875.3335
$ws.Range("J84").Value = 250025000
$ws.Range("K84").Value = 52878.0015
$ws.Range("L84").Value = 2250225000
$ws.Range("M84").Value = -47262.0015
$ws.Range("N84").Value = -2250236232
$ws.Range("H122").Value = 660
$ws.Range("I122").Value = 579.8889
$ws.Range("J122").Value = 840.25
$ws.Range("K122").Value = 5219.0001
$ws.Range("L122").Value = 7562.25
$ws.Range("M122").Value = -2769.0001
$ws.Range("N122").Value = -12462.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7165.3335
$ws.Range("I70").Value = 6998
$ws.Range("J70").Value = 7500
$ws.Range("K70").Value = 6998
$ws.Range("L70").Value = 7500
$ws.Range("M70").Value = -6728
$ws.Range("N70").Value = -8040
$ws.Range("H73").Value = 7165.3335
$ws.Range("I73").Value = 6998
$ws.Range("J73").Value = 7500
$ws.Range("K73").Value = 6998
$ws.Range("L73").Value = 7500
$ws.Range("M73").Value = -6062
$ws.Range("N73").Value = -9372
$ws.Range("H113").Value = 2212.1333
$ws.Range("I113").Value = 2022.1111
$ws.Range("K113").Value = 2022.1111
$ws.Range("M113").Value = 147.8888999999999
$ws.Range("H122").Value = 1892.4546
$ws.Range("I122").Value = 1651.6666
$ws.Range("J122").Value = 2181.4
$ws.Range("K122").Value = 4954.9998
$ws.Range("L122").Value = 6544.200000000001
$ws.Range("M122").Value = -2504.9998
$ws.Range("N122").Value = -11444.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("M16").Value = -830
$ws.Range("H39").Value = 10000
$ws.Range("J39").Value = 10000
$ws.Range("L39").Value = 10000
$ws.Range("N39").Value = -10920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 8697.5
$ws.Range("I14").Value = 1400
$ws.Range("K14").Value = 1400
$ws.Range("M14").Value = -1232
$ws.Range("H126").Value = 1303.6666
$ws.Range("I126").Value = 1215.3
$ws.Range("J126").Value = 1745.5
$ws.Range("K126").Value = 3645.9
$ws.Range("L126").Value = 5236.5
$ws.Range("M126").Value = -1175.9
$ws.Range("N126").Value = -10176.5
$ws.Range("H132").Value = 4620.976
$ws.Range("I132").Value = 5623.517
$ws.Range("J132").Value = 2384.5386
$ws.Range("K132").Value = 16870.551
$ws.Range("L132").Value = 7153.6158
$ws.Range("M132").Value = -14340.551
$ws.Range("N132").Value = -12213.6158
